$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '27.899.77'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -0.22%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.633.25'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -0.31%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '211.51'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.39%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.520'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.84%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '23.44'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +0.49%  '
$ws.Range('E9').Value = '  -0.86%  '
$ws.Range('E10').Value = '  -0.19%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0883'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +0.32%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '1.864.98'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -0.32%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '1.641.86'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -0.74%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.565'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -0.74%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '65.40'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +0.06%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '27.902.91'
$c.Style = 'Normal'
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '229.06'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.99%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '7.69'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +1.95%  '
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('E22').Value = '  -0.88%  '
$ws.Range('E23').Value = '  -3.37%  '
$ws.Range('E24').Value = '  -0.34%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '155.55'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +1.28%  '
$ws.Range('E26').Value = '  -1.36%  '
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('E28').Value = '  -0.53%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('E30').Value = '  -0.38%  '
$ws.Range('E31').Value = '  -0.40%  '
$ws.Range('E32').Value = '  +0.77%  '
$ws.Range('E33').Value = '  +1.12%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.393.87'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -0.98%  '
$ws.Range('E35').Value = '  +0.94%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.04'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +11.90%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '2.35'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -0.46%  '
$ws.Range('E38').Value = '  +1.04%  '
$ws.Range('E39').Value = '  -0.51%  '
$ws.Range('E40').Value = '  -3.13%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.04%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '1.01'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -0.61%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '65.89'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -1.74%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '1.83'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +0.88%  '
$ws.Range('E45').Value = '  -1.53%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '1.774.36'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -0.36%  '
$ws.Range('E47').Value = '  -3.06%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '88.86'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +0.99%  '
$ws.Range('E49').Value = '  +1.56%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '7.68'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +1.44%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.0504'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -0.24%  '
